# Applies market-price / profit-column updates captured from the scheduled
# Zalera market-data refresh. Values are plain numbers (no formulas in this
# workbook), so each touched cell is written directly via .Value; the one
# cell the refresh dropped (no longer computable) is cleared instead.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 5033.143
$ws.Range("I48").Value = 5077.3335
$ws.Range("J48").Value = 5000
$ws.Range("K48").Value = 15232.0005
$ws.Range("L48").Value = 15000
$ws.Range("M48").Value = -14940.0005
$ws.Range("N48").Value = -15584
$ws.Range("H56").Value = 5033.143
$ws.Range("I56").Value = 5077.3335
$ws.Range("J56").Value = 5000
$ws.Range("K56").Value = 15232.0005
$ws.Range("L56").Value = 15000
$ws.Range("M56").Value = -14698.0005
$ws.Range("N56").Value = -16068
$ws.Range("H64").Value = 6336.143
$ws.Range("I64").Value = 6110.8
$ws.Range("J64").Value = 6899.5
$ws.Range("K64").Value = 6110.8
$ws.Range("L64").Value = 6899.5
$ws.Range("M64").Value = -5862.8
$ws.Range("N64").Value = -7395.5
$ws.Range("H67").Value = 6336.143
$ws.Range("I67").Value = 6110.8
$ws.Range("J67").Value = 6899.5
$ws.Range("K67").Value = 6110.8
$ws.Range("L67").Value = 6899.5
$ws.Range("M67").Value = -5252.8
$ws.Range("N67").Value = -8615.5
$ws.Range("H70").Value = 6442.25
$ws.Range("I70").Value = 9000
$ws.Range("K70").Value = 27000
$ws.Range("M70").Value = -26730
$ws.Range("H73").Value = 6442.25
$ws.Range("I73").Value = 9000
$ws.Range("K73").Value = 27000
$ws.Range("M73").Value = -26064
$ws.Range("H74").Value = 10374.75
$ws.Range("I74").Value = 10749.5
$ws.Range("K74").Value = 10749.5
$ws.Range("M74").Value = -9813.5
$ws.Range("I76").Value = 20005000
$ws.Range("J76").Value = 5555
$ws.Range("K76").Value = 20005000
$ws.Range("L76").Value = 5555
$ws.Range("M76").Value = -20004685
$ws.Range("N76").Value = -6185
$ws.Range("H77").Value = 10374.75
$ws.Range("I77").Value = 10749.5
$ws.Range("K77").Value = 53747.5
$ws.Range("M77").Value = -49067.5
$ws.Range("I79").Value = 20005000
$ws.Range("J79").Value = 5555
$ws.Range("K79").Value = 20005000
$ws.Range("L79").Value = 5555
$ws.Range("M79").Value = -20003908
$ws.Range("N79").Value = -7739
$ws.Range("H103").Value = 1450.75
$ws.Range("I103").Value = 1259.4286
$ws.Range("K103").Value = 3778.2858
$ws.Range("M103").Value = -3192.2858
$ws.Range("H112").Value = 2846.6296
$ws.Range("J112").Value = 3011.2083
$ws.Range("L112").Value = 9033.624899999999
$ws.Range("N112").Value = -11249.6249
$ws.Range("H116").Value = 4577.3335
$ws.Range("I116").Value = 4240
$ws.Range("J116").Value = 4999
$ws.Range("K116").Value = 4240
$ws.Range("L116").Value = 4999
$ws.Range("M116").Value = -798
$ws.Range("N116").Value = -11883
$ws.Range("H137").Value = 17867914
$ws.Range("I137").Value = 125001250
$ws.Range("J137").Value = 12358.167
$ws.Range("K137").Value = 375003750
$ws.Range("L137").Value = 37074.501
$ws.Range("M137").Value = -375001200
$ws.Range("N137").Value = -42174.501

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 32880.656
$ws.Range("I32").Value = 32880.656
$ws.Range("K32").Value = 32880.656
$ws.Range("M32").Value = -32593.656
$ws.Range("H63").Value = 8283.666999999999
$ws.Range("I63").Value = 3149.6667
$ws.Range("J63").Value = 9995
$ws.Range("K63").Value = 3149.6667
$ws.Range("L63").Value = 9995
$ws.Range("M63").Value = -2463.6667
$ws.Range("N63").Value = -11367
$ws.Range("H66").Value = 8283.666999999999
$ws.Range("I66").Value = 3149.6667
$ws.Range("J66").Value = 9995
$ws.Range("K66").Value = 15748.3335
$ws.Range("L66").Value = 49975
$ws.Range("M66").Value = -12316.3335
$ws.Range("N66").Value = -56839
$ws.Range("H88").Value = 6207.778
$ws.Range("J88").Value = 9280.799999999999
$ws.Range("L88").Value = 9280.799999999999
$ws.Range("N88").Value = -10092.8
$ws.Range("H91").Value = 6207.778
$ws.Range("J91").Value = 9280.799999999999
$ws.Range("L91").Value = 9280.799999999999
$ws.Range("N91").Value = -12088.8
$ws.Range("H96").Value = 28518
$ws.Range("J96").Value = 33428.8
$ws.Range("L96").Value = 33428.8
$ws.Range("N96").Value = -38920.8
$ws.Range("H132").Value = 3378.3
$ws.Range("I132").Value = 2307.56
$ws.Range("K132").Value = 6922.68
$ws.Range("M132").Value = -4392.68

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 36617
$ws.Range("I26").Value = 36617
$ws.Range("K26").Value = 36617
$ws.Range("M26").Value = -36325
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H86").Value = 186063.36
$ws.Range("I86").Value = 4669.7
$ws.Range("J86").Value = 2000000
$ws.Range("K86").Value = 4669.7
$ws.Range("L86").Value = 2000000
$ws.Range("M86").Value = -3546.7
$ws.Range("N86").Value = -2002246
$ws.Range("H89").Value = 186063.36
$ws.Range("I89").Value = 4669.7
$ws.Range("J89").Value = 2000000
$ws.Range("K89").Value = 23348.5
$ws.Range("L89").Value = 10000000
$ws.Range("M89").Value = -17732.5
$ws.Range("N89").Value = -10011232
$ws.Range("H95").Value = 19393.25
$ws.Range("J95").Value = 19393.25
$ws.Range("L95").Value = 19393.25
$ws.Range("N95").Value = -24885.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4232.1816
$ws.Range("I16").Value = 3945.4
$ws.Range("K16").Value = 3945.4
$ws.Range("M16").Value = -3658.4
$ws.Range("H58").Value = 4302.3
$ws.Range("I58").Value = 2696.6155
$ws.Range("J58").Value = 7284.2856
$ws.Range("K58").Value = 2696.6155
$ws.Range("L58").Value = 7284.2856
$ws.Range("M58").Value = -2493.6155
$ws.Range("N58").Value = -7690.2856
$ws.Range("H96").Value = 34927.11
$ws.Range("J96").Value = 34927.11
$ws.Range("L96").Value = 34927.11
$ws.Range("N96").Value = -40419.11
$ws.Range("H99").Value = 4999.5
$ws.Range("I99").Value = 4999
$ws.Range("K99").Value = 4999
$ws.Range("M99").Value = -3501
$ws.Range("H113").Value = 4232.1816
$ws.Range("I113").Value = 3945.4
$ws.Range("K113").Value = 3945.4
$ws.Range("M113").Value = -1775.4
$ws.Range("H126").Value = 4999.5
$ws.Range("I126").Value = 4999
$ws.Range("K126").Value = 14997
$ws.Range("M126").Value = -12527
$ws.Range("H136").Value = 4302.3
$ws.Range("I136").Value = 2696.6155
$ws.Range("J136").Value = 7284.2856
$ws.Range("K136").Value = 8089.8465
$ws.Range("L136").Value = 21852.8568
$ws.Range("M136").Value = -5539.8465
$ws.Range("N136").Value = -26952.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 18306.234
$ws.Range("I26").Value = 142.76923
$ws.Range("K26").Value = 428.30769
$ws.Range("M26").Value = -140.30769
$ws.Range("H92").Value = 1327.1428
$ws.Range("J92").Value = 963.8333
$ws.Range("L92").Value = 2891.4999
$ws.Range("N92").Value = -5387.4999
$ws.Range("H109").Value = 3694.5833
$ws.Range("I109").Value = 2633.5
$ws.Range("K109").Value = 7900.5
$ws.Range("M109").Value = -6860.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 41681204
$ws.Range("I40").Value = 41681204
$ws.Range("K40").Value = 41681204
$ws.Range("M40").Value = -41681068
$ws.Range("H61").Value = 2382.111
$ws.Range("J61").Value = 4902.5
$ws.Range("L61").Value = 4902.5
$ws.Range("N61").Value = -5306.5
$ws.Range("H74").Value = 50841
$ws.Range("I74").Value = 47914
$ws.Range("K74").Value = 47914
$ws.Range("M74").Value = -46916
$ws.Range("H77").Value = 50841
$ws.Range("I77").Value = 47914
$ws.Range("K77").Value = 143742
$ws.Range("M77").Value = -138750
$ws.Range("H82").Value = 2378.2
$ws.Range("J82").Value = 2975
$ws.Range("L82").Value = 2975
$ws.Range("N82").Value = -3697
$ws.Range("H85").Value = 2378.2
$ws.Range("J85").Value = 2975
$ws.Range("L85").Value = 2975
$ws.Range("N85").Value = -5471
$ws.Range("H113").Value = 2382.111
$ws.Range("J113").Value = 4902.5
$ws.Range("L113").Value = 4902.5
$ws.Range("N113").Value = -9242.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 28844.54
$ws.Range("I63").Value = 25555.445
$ws.Range("J63").Value = 36245
$ws.Range("K63").Value = 25555.445
$ws.Range("L63").Value = 36245
$ws.Range("M63").Value = -24931.445
$ws.Range("N63").Value = -37493
$ws.Range("H66").Value = 28844.54
$ws.Range("I66").Value = 25555.445
$ws.Range("J66").Value = 36245
$ws.Range("K66").Value = 76666.33499999999
$ws.Range("L66").Value = 108735
$ws.Range("M66").Value = -73546.33499999999
$ws.Range("N66").Value = -114975
$ws.Range("H81").Value = 5624.2144
$ws.Range("I81").Value = 3088.9
$ws.Range("J81").Value = 11962.5
$ws.Range("K81").Value = 6177.8
$ws.Range("L81").Value = 23925
$ws.Range("M81").Value = -5116.8
$ws.Range("N81").Value = -26047
$ws.Range("H84").Value = 5624.2144
$ws.Range("I84").Value = 3088.9
$ws.Range("J84").Value = 11962.5
$ws.Range("K84").Value = 30889
$ws.Range("L84").Value = 119625
$ws.Range("M84").Value = -25585
$ws.Range("N84").Value = -130233
$ws.Range("H86").Value = 77777
$ws.Range("J86").Value = 77777
$ws.Range("L86").Value = 77777
$ws.Range("N86").Value = -80023
$ws.Range("H89").Value = 77777
$ws.Range("J89").Value = 77777
$ws.Range("L89").Value = 388885
$ws.Range("N89").Value = -400117
$ws.Range("H100").Value = 1270.2
$ws.Range("I100").Value = 1189.1111
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 2378.2222
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -1837.2222
$ws.Range("N100").Value = -5082

